# Fruta / hortaliza, semanal
# Insert 3 new weekly rows (variety "Artic Star") above the existing row 544,
# pushing the previously-existing rows 544:561 down to 547:564.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 544-546 (existing data shifts down to 547-564).
$ws.Range("A544:A546").EntireRow.Insert()

# New row 544 - Artic Star / Especial
$ws.Cells.Item(544, 1).Value = 8
$ws.Cells.Item(544, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(544, 3).Value = "Coquimbo"
$ws.Cells.Item(544, 4).Value = 45267
$ws.Cells.Item(544, 5).Value = 4
$ws.Cells.Item(544, 6).Value = "Fruta"
$ws.Cells.Item(544, 7).Value = 100103
$ws.Cells.Item(544, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(544, 9).Value = 100103006
$ws.Cells.Item(544, 10).Value = "Nectarín"
$ws.Cells.Item(544, 11).Value = "Artic Star"
$ws.Cells.Item(544, 12).Value = "Especial"
$ws.Cells.Item(544, 13).Value = 16
$ws.Cells.Item(544, 14).Value = 400000
$ws.Cells.Item(544, 15).Value = 410000
$ws.Cells.Item(544, 16).Value = 405000
$ws.Cells.Item(544, 17).Value = "`$/bins (420 kilos)"
$ws.Cells.Item(544, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(544, 19).Value = 964
$ws.Cells.Item(544, 20).Value = 420

# New row 545 - Artic Star / Primera
$ws.Cells.Item(545, 1).Value = 8
$ws.Cells.Item(545, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(545, 3).Value = "Coquimbo"
$ws.Cells.Item(545, 4).Value = 45267
$ws.Cells.Item(545, 5).Value = 4
$ws.Cells.Item(545, 6).Value = "Fruta"
$ws.Cells.Item(545, 7).Value = 100103
$ws.Cells.Item(545, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(545, 9).Value = 100103006
$ws.Cells.Item(545, 10).Value = "Nectarín"
$ws.Cells.Item(545, 11).Value = "Artic Star"
$ws.Cells.Item(545, 12).Value = "Primera"
$ws.Cells.Item(545, 13).Value = 10
$ws.Cells.Item(545, 14).Value = 370000
$ws.Cells.Item(545, 15).Value = 380000
$ws.Cells.Item(545, 16).Value = 375000
$ws.Cells.Item(545, 17).Value = "`$/bins (420 kilos)"
$ws.Cells.Item(545, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(545, 19).Value = 893
$ws.Cells.Item(545, 20).Value = 420

# New row 546 - Artic Star / Segunda
$ws.Cells.Item(546, 1).Value = 8
$ws.Cells.Item(546, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(546, 3).Value = "Coquimbo"
$ws.Cells.Item(546, 4).Value = 45267
$ws.Cells.Item(546, 5).Value = 4
$ws.Cells.Item(546, 6).Value = "Fruta"
$ws.Cells.Item(546, 7).Value = 100103
$ws.Cells.Item(546, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(546, 9).Value = 100103006
$ws.Cells.Item(546, 10).Value = "Nectarín"
$ws.Cells.Item(546, 11).Value = "Artic Star"
$ws.Cells.Item(546, 12).Value = "Segunda"
$ws.Cells.Item(546, 13).Value = 16
$ws.Cells.Item(546, 14).Value = 340000
$ws.Cells.Item(546, 15).Value = 350000
$ws.Cells.Item(546, 16).Value = 345000
$ws.Cells.Item(546, 17).Value = "`$/bins (420 kilos)"
$ws.Cells.Item(546, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(546, 19).Value = 821
$ws.Cells.Item(546, 20).Value = 420

Write-Output "ok"
